$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reserve currently-untranslated files for translator "cfarl" (fill column D - Tradutor)
$ws.Range("D561:D582").Value = "cfarl"
$ws.Range("D587:D594").Value = "cfarl"
$ws.Range("D596:D606").Value = "cfarl"

# Mark files mes_inn01_us.u16 .. mes_inn04_us.u16 (rows 755-758) as translated ("Sim") by cfarl
$ws.Range("C755:C758").Value = "Sim"
$ws.Range("D755:D758").Value = "cfarl"

# Update selection to reflect where the editor left off
$ws.Range("E758").Select() | Out-Null
